$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Metadata" sheet: bump the Date value and update the Contact URL
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value  = "2025-05-12T20:06:16+02:00"
$meta.Range("B10").Value = "Constantin Renner (http://example.org/example-publisher, constantinrenner1@gmail.com)"

# ---------------------------------------------------------------------------
# 2) "Include #0" sheet: add three new SNOMED concept rows (Slow / Fast /
#    Normal) right after the existing "Absent (qualifier value)" row, ahead
#    of the trailing blank separator row + "System URI" row.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Include #0")

# Push the existing row 5 (blank separator) and row 6 (System URI row) down
# by three rows so we end up with free rows 5-7 to fill in. Using Insert()
# (rather than overwriting) keeps those two rows' original cell encoding
# completely untouched.
$ws.Range("A5:B7").Insert()

# Clone the formatting of an existing data row (row 4) onto the freshly
# inserted rows 5-7 so they pick up the same style (borders / wrap text)
# as every other concept row instead of a blank default style.
$ws.Range("A4:B4").Copy($ws.Range("A5:B7"))

# Helper cell used to stage text so numeric-looking codes ("255361000", ...)
# are written as TEXT (matching how the other SNOMED codes in this sheet are
# stored) instead of being auto-converted to numbers, while not disturbing
# the row/column formatting already in place on A5:B7.
$stage = $ws.Range("Z1")
$stage.NumberFormat = "@"

$stage.Value = "255361000"
$stage.Copy()
$ws.Range("A5").PasteSpecial(-4163)   # xlPasteValues

$stage.Value = "277748003"
$stage.Copy()
$ws.Range("A6").PasteSpecial(-4163)

$stage.Value = "17621005"
$stage.Copy()
$ws.Range("A7").PasteSpecial(-4163)

$stage.Clear()

$ws.Range("B5").Value = "Slow (qualifier value)"
$ws.Range("B6").Value = "Fast (qualifier value)"
$ws.Range("B7").Value = "  Normal (qualifier value)"
